# Change_Log.xlsx edit: document the BEQ/branch-instruction progress entry
# (adds a Notes cell for the previous entry + a new Changes/Notes row pair).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Change Log")

# --- Row 20: fill in the Notes cell that was left as a placeholder ---
$ws.Range("D20").Value = "Notes`n- Stopped at about min 27 of vid 20. My implementation uses an extra cycle than needed, which I found out from the video pretty much`nBugs`n- "
$ws.Rows.Item(20).RowHeight = 86.4

# --- Row 21: turn the next blank template row into the new entry ---
$ws.Range("A21").Value = 45958
$ws.Range("B21").Value = "Changes`n- MODIFIED: MyMIF.mif, alu.vhd, alu_control.vhd, Controller.vhd         `n- COMPLETED: Implemented all the branch instructions (BEQ, BNE, BLEZ, BGEZ, BLTZ, BGTZ). Only tested BEQ in the .mif file so far                                                                                                                                                                                                                              "
$ws.Range("D21").Value = "Notes`n- I have only tested BEQ so far, but all of the branch instructions are very similar so I expect the rest of them to work as well. After the initial instruction decode state in the controller, I was able to implement with only 1 other branch state after initially making different states for each instruction`n- I mean tbh the only thing that's really left is to test the rest of the branch instructions, which I highly expect to work. So really I'm pretty much done lol. YAYYYYY!!!`nBugs`n- "
$ws.Rows.Item(21).RowHeight = 144

# --- Update the active selection to the next Notes placeholder ---
$ws.Range("D22").Select()
